$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "2025-04-28 23:44:21"
$ws.Range("B34").Value = 161
